$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task text edits ---
# Row 26: "Change the schematic symbol for devices so gnd is on the bottom"
#         -> "Change the schematic symbol for devices so gnd isn't on the bottom"
$ws.Cells.Item(26, 1).Value = "Change the schematic symbol for devices so gnd isn't on the bottom"

# Row 28: "Delete test point nets, replace with TP symbols" -> "TP symbols"
$ws.Cells.Item(28, 1).Value = "TP symbols"

# --- Mark additional tasks as done (Status column B: 0 -> 1) ---
$doneRows = @(22, 23, 25, 27, 30, 33, 34, 35, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 49, 50, 51, 56)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 2).Value = 1
}

# --- Update the view state to match where the author left the cursor/scroll ---
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("D56").Select()
